# Generate Report for Handoff
# A new handoff round produced a fresh report GUID / content hash; refresh
# the Overview + per-locale status sheets (file names, timestamps and the
# hyperlinks that display those file names) to point at the new report.

$wb = $excel.ActiveWorkbook

$newGuid = "b48c9462-c430-49ae-95b3-2a188527c579"
$newContentHash = "2fa471eccc545a3cb8a30ce65877f17157450074"

$newHandoffDate = "2016-03-23 17:10:38"   # Overview!D2 and de-de!E2
$newZhDate = "2016-03-23 17:10:30"        # zh-cn!E2

# Existing external hyperlink targets (unchanged by this edit) - reapplied
# after rewriting each cell so the link keeps pointing at the same URL while
# only its visible display text is refreshed to the new file name.
$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/c04d5a9ca12eb81ce35c2d00359cc36dc9751efc/e2e/576178ce-7fb2-4e34-ac16-ec2cb9fe4a70.md"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bdfba8dfad230577c2963fa321eb840bdfb20651/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/576178ce-7fb2-4e34-ac16-ec2cb9fe4a70.ed04ee567b96368f8b5556ae18779750b06aedca.zh-cn.xlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/281076494b011462186df01e660f3c3a139ec82e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/576178ce-7fb2-4e34-ac16-ec2cb9fe4a70.ed04ee567b96368f8b5556ae18779750b06aedca.de-de.xlf"

function Update-HyperlinkCell($range, $target, $display) {
    # Replacing a hyperlink's Address/TextToDisplay in place isn't reliable
    # through this host, so drop the old link, write the new display text as
    # the cell value, then re-attach a hyperlink with the original target.
    $range.Hyperlinks.Delete()
    $range.Value = $display
    $h = $range.Hyperlinks.Item(1)
    $h.Address = $target
    $h.TextToDisplay = $display
}

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
Update-HyperlinkCell $wsOverview.Range("A2") $mdTarget "$newGuid.md"
$wsOverview.Range("D2").Value = $newHandoffDate

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HyperlinkCell $wsZhCn.Range("A2") $mdTarget "$newGuid.md"
Update-HyperlinkCell $wsZhCn.Range("D2") $zhXlfTarget "$newGuid.$newContentHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = $newZhDate

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HyperlinkCell $wsDeDe.Range("A2") $mdTarget "$newGuid.md"
Update-HyperlinkCell $wsDeDe.Range("D2") $deXlfTarget "$newGuid.$newContentHash.de-de.xlf"
$wsDeDe.Range("E2").Value = $newHandoffDate
